# Commit: "replace trelp with skay"
#
# Visible change: the "Trelp" title text box on slides 1-7 and 9 becomes
# "Skay". (The cached "datetimeFigureOut" fields on the slide master and
# every slide layout also shifted from 2/23/16 to 5/25/16 - a side effect
# of PowerPoint refreshing the auto-date fields on save - so we bring
# those along too.)

$p = $ppt.ActivePresentation

# --- 1. Trelp -> Skay on every slide that has it -----------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "Trelp") {
                $shp.TextFrame.TextRange.Text = "Skay"
            }
        }
    }
}

# --- 2. Refresh the cached date field (2/23/16 -> 5/25/16) -------------
# ppPlaceholderDate = 16
function Update-DatePlaceholder($shapes, $newText) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
        }
        if ($isDatePlaceholder -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "2/23/16") {
                $shp.TextFrame.TextRange.Text = "5/25/16"
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes "5/25/16"

# Every slide layout (go through Presentation.SlideMaster.CustomLayouts
# so each layout is addressed individually)
for ($j = 1; $j -le $p.SlideMaster.CustomLayouts.Count; $j++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($j)
    Update-DatePlaceholder $layout.Shapes "5/25/16"
}
